$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 74899.08
$ws.Range("I80").Value = 17389.143
$ws.Range("J80").Value = 141994
$ws.Range("K80").Value = 52167.429
$ws.Range("L80").Value = 425982
$ws.Range("M80").Value = -51169.429
$ws.Range("N80").Value = -427978

$ws.Range("H83").Value = 74899.08
$ws.Range("I83").Value = 17389.143
$ws.Range("J83").Value = 141994
$ws.Range("K83").Value = 156502.287
$ws.Range("L83").Value = 1277946
$ws.Range("M83").Value = -151510.287
$ws.Range("N83").Value = -1287930

$ws.Range("H86").Value = 198416540
$ws.Range("I86").Value = 333334300
$ws.Range("K86").Value = 333334300
$ws.Range("M86").Value = -333333177

$ws.Range("H89").Value = 198416540
$ws.Range("I89").Value = 333334300
$ws.Range("K89").Value = 1666671500
$ws.Range("M89").Value = -1666665884

$ws.Range("H98").Value = 5576
$ws.Range("J98").Value = 500
$ws.Range("L98").Value = 500
$ws.Range("N98").Value = -3496

$ws.Range("H107").Value = 17859776
$ws.Range("I107").Value = 7355283
$ws.Range("J107").Value = 62503876
$ws.Range("K107").Value = 7355283
$ws.Range("L107").Value = 62503876
$ws.Range("M107").Value = -7353363
$ws.Range("N107").Value = -62507716

$ws.Range("H111").Value = 10417835
$ws.Range("J111").Value = 1099.75
$ws.Range("L111").Value = 3299.25
$ws.Range("N111").Value = -9433.25

$ws.Range("H112").Value = 9079.583000000001
$ws.Range("J112").Value = 9693.409
$ws.Range("L112").Value = 29080.227
$ws.Range("N112").Value = -31296.227

$ws.Range("H122").Value = 5576
$ws.Range("J122").Value = 500
$ws.Range("L122").Value = 1500
$ws.Range("N122").Value = -6400

$ws.Range("H132").Value = 1395.6735
$ws.Range("I132").Value = 1418.3617
$ws.Range("K132").Value = 4255.0851
$ws.Range("M132").Value = -1725.0851

$ws.Range("H137").Value = 1073.4
$ws.Range("I137").Value = 1087.1111
$ws.Range("J137").Value = 950
$ws.Range("K137").Value = 3261.3333
$ws.Range("L137").Value = 2850
$ws.Range("M137").Value = -711.3333000000002
$ws.Range("N137").Value = -7950

$ws.Range("H138").Value = 3458682.8
$ws.Range("J138").Value = 3714572.2
$ws.Range("L138").Value = 11143716.6
$ws.Range("N138").Value = -11153996.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2162
$ws.Range("J63").Value = 2616
$ws.Range("L63").Value = 2616
$ws.Range("N63").Value = -3988

$ws.Range("H66").Value = 2162
$ws.Range("J66").Value = 2616
$ws.Range("L66").Value = 13080
$ws.Range("N66").Value = -19944

$ws.Range("H74").Value = 37874.965
$ws.Range("I74").Value = 51831.55
$ws.Range("J74").Value = 2983.5
$ws.Range("K74").Value = 51831.55
$ws.Range("L74").Value = 2983.5
$ws.Range("M74").Value = -50957.55
$ws.Range("N74").Value = -4731.5

$ws.Range("H77").Value = 37874.965
$ws.Range("I77").Value = 51831.55
$ws.Range("J77").Value = 2983.5
$ws.Range("K77").Value = 259157.75
$ws.Range("L77").Value = 14917.5
$ws.Range("M77").Value = -254789.75
$ws.Range("N77").Value = -23653.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 55397.43
$ws.Range("I86").Value = 79464.71000000001
$ws.Range("J86").Value = 7262.857
$ws.Range("K86").Value = 79464.71000000001
$ws.Range("L86").Value = 7262.857
$ws.Range("M86").Value = -78341.71000000001
$ws.Range("N86").Value = -9508.857

$ws.Range("H89").Value = 55397.43
$ws.Range("I89").Value = 79464.71000000001
$ws.Range("J89").Value = 7262.857
$ws.Range("K89").Value = 397323.55
$ws.Range("L89").Value = 36314.285
$ws.Range("M89").Value = -391707.55
$ws.Range("N89").Value = -47546.285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 11905132
$ws.Range("I22").Value = 271.22223
$ws.Range("J22").Value = 47619710
$ws.Range("K22").Value = 271.22223
$ws.Range("L22").Value = 47619710
$ws.Range("M22").Value = 78.77776999999998
$ws.Range("N22").Value = -47620410

$ws.Range("H31").Value = 4886.3955
$ws.Range("I31").Value = 1017.4211
$ws.Range("J31").Value = 7949.3335
$ws.Range("K31").Value = 1017.4211
$ws.Range("L31").Value = 7949.3335
$ws.Range("M31").Value = -722.4211
$ws.Range("N31").Value = -8539.333500000001

$ws.Range("H34").Value = 4886.3955
$ws.Range("I34").Value = 1017.4211
$ws.Range("J34").Value = 7949.3335
$ws.Range("K34").Value = 1017.4211
$ws.Range("L34").Value = 7949.3335
$ws.Range("M34").Value = -815.4211
$ws.Range("N34").Value = -8353.333500000001

$ws.Range("H62").Value = 1995.75
$ws.Range("I62").Value = 1995.75
$ws.Range("K62").Value = 1995.75
$ws.Range("M62").Value = -1371.75

$ws.Range("H65").Value = 1995.75
$ws.Range("I65").Value = 1995.75
$ws.Range("K65").Value = 9978.75
$ws.Range("M65").Value = -6858.75

$ws.Range("H107").Value = 958.12823
$ws.Range("I107").Value = 302.78262
$ws.Range("J107").Value = 1900.1875
$ws.Range("K107").Value = 302.78262
$ws.Range("L107").Value = 1900.1875
$ws.Range("M107").Value = 1617.21738
$ws.Range("N107").Value = -5740.1875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2004.0714
$ws.Range("J68").Value = 2006.8889
$ws.Range("L68").Value = 6020.6667
$ws.Range("N68").Value = -7642.6667

$ws.Range("H71").Value = 2004.0714
$ws.Range("J71").Value = 2006.8889
$ws.Range("L71").Value = 18062.0001
$ws.Range("N71").Value = -26174.0001

$ws.Range("H113").Value = 1614.1177
$ws.Range("I113").Value = 919.55554
$ws.Range("J113").Value = 2395.5
$ws.Range("K113").Value = 2758.66662
$ws.Range("L113").Value = 7186.5
$ws.Range("M113").Value = -588.66662
$ws.Range("N113").Value = -11526.5

$ws.Range("H129").Value = 63592.938
$ws.Range("I129").Value = 513.2222
$ws.Range("J129").Value = 144695.42
$ws.Range("K129").Value = 1539.6666
$ws.Range("L129").Value = 434086.26
$ws.Range("M129").Value = 3460.3334
$ws.Range("N129").Value = -444086.26

$ws.Range("H131").Value = 2010.46
$ws.Range("J131").Value = 2156.4595
$ws.Range("L131").Value = 6469.3785
$ws.Range("N131").Value = -16549.3785

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 23487.215
$ws.Range("J63").Value = 8499
$ws.Range("L63").Value = 8499
$ws.Range("N63").Value = -9871

$ws.Range("H66").Value = 23487.215
$ws.Range("J66").Value = 8499
$ws.Range("L66").Value = 25497
$ws.Range("N66").Value = -32361

$ws.Range("H80").Value = 2893.4443
$ws.Range("I80").Value = 2468.1428
$ws.Range("J80").Value = 3164.0908
$ws.Range("K80").Value = 2468.1428
$ws.Range("L80").Value = 3164.0908
$ws.Range("M80").Value = -1470.1428
$ws.Range("N80").Value = -5160.0908

$ws.Range("H83").Value = 2893.4443
$ws.Range("I83").Value = 2468.1428
$ws.Range("J83").Value = 3164.0908
$ws.Range("K83").Value = 12340.714
$ws.Range("L83").Value = 15820.454
$ws.Range("M83").Value = -7348.714
$ws.Range("N83").Value = -25804.454

$ws.Range("H113").Value = 5502.1
$ws.Range("I113").Value = 5002.3335
$ws.Range("K113").Value = 5002.3335
$ws.Range("M113").Value = -2832.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6459.303
$ws.Range("J122").Value = 6714.857
$ws.Range("L122").Value = 20144.571
$ws.Range("N122").Value = -25044.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 3000
$ws.Range("I8").Value = 5000
$ws.Range("J8").Value = 1000
$ws.Range("K8").Value = 5000
$ws.Range("L8").Value = 1000
$ws.Range("M8").Value = -4860
$ws.Range("N8").Value = -1280

$ws.Range("H15").Value = 24599.25
$ws.Range("I15").Value = 24599.25
$ws.Range("K15").Value = 24599.25
$ws.Range("M15").Value = -24311.25

$ws.Range("H113").Value = 19200.572
$ws.Range("J113").Value = 1579.2222
$ws.Range("L113").Value = 4737.6666
$ws.Range("N113").Value = -9077.6666

$ws.Range("H122").Value = 81899.96000000001
$ws.Range("J122").Value = 15636.363
$ws.Range("L122").Value = 46909.089
$ws.Range("N122").Value = -51809.089
